$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds quarter-start dates stored as literal text (e.g. "01/07/2013")
# rather than real Excel dates. A leading apostrophe forces the new value to be
# entered as literal text too, instead of being auto-converted to a date serial.

$ws.Range("C2").Value = "'01/10/2013"
$ws.Range("D2").Value = 78.75336125633515
$ws.Range("C3").Value = "'01/10/2015"
$ws.Range("D3").Value = 79.45616811315006
$ws.Range("C4").Value = "'01/10/2017"
$ws.Range("D4").Value = 79.99796933698853
$ws.Range("C5").Value = "'01/10/2019"
$ws.Range("D5").Value = 80.4452653074825
$ws.Range("C6").Value = "'01/10/2021"
$ws.Range("D6").Value = 80.81270990862525
$ws.Range("C7").Value = "'01/10/2023"
$ws.Range("D7").Value = 81.13527944735173
$ws.Range("C8").Value = "'01/10/2013"
$ws.Range("D8").Value = 49.14231734258324
$ws.Range("C9").Value = "'01/10/2015"
$ws.Range("D9").Value = 49.92480907402129
$ws.Range("C10").Value = "'01/10/2017"
$ws.Range("D10").Value = 50.61282508739103
$ws.Range("C11").Value = "'01/10/2019"
$ws.Range("D11").Value = 51.13267992212379
$ws.Range("C12").Value = "'01/10/2021"
$ws.Range("D12").Value = 50.54599696042929
$ws.Range("C13").Value = "'01/10/2023"
$ws.Range("D13").Value = 50.46478162897979
$ws.Range("C14").Value = "'01/10/2013"
$ws.Range("D14").Value = 46.06795485670302
$ws.Range("C15").Value = "'01/10/2015"
$ws.Range("D15").Value = 45.39271286895155
$ws.Range("C16").Value = "'01/10/2017"
$ws.Range("D16").Value = 44.59142576717965
$ws.Range("C17").Value = "'01/10/2019"
$ws.Range("D17").Value = 45.46666222385126
$ws.Range("C18").Value = "'01/10/2021"
$ws.Range("D18").Value = 44.91200255173838
$ws.Range("C19").Value = "'01/10/2023"
$ws.Range("D19").Value = 46.72570712048232
$ws.Range("C20").Value = "'01/10/2013"
$ws.Range("D20").Value = 3.074362485880224
$ws.Range("C21").Value = "'01/10/2015"
$ws.Range("D21").Value = 4.532096205069736
$ws.Range("C22").Value = "'01/10/2017"
$ws.Range("D22").Value = 6.020915829018175
$ws.Range("C23").Value = "'01/10/2019"
$ws.Range("D23").Value = 5.666017698272538
$ws.Range("C24").Value = "'01/10/2021"
$ws.Range("D24").Value = 5.633994408690921
$ws.Range("C25").Value = "'01/10/2023"
$ws.Range("D25").Value = 3.73953720797879
$ws.Range("C26").Value = "'01/10/2013"
$ws.Range("D26").Value = 29.61104391375191
$ws.Range("C27").Value = "'01/10/2015"
$ws.Range("D27").Value = 29.53185048308941
$ws.Range("C28").Value = "'01/10/2017"
$ws.Range("D28").Value = 29.3851442495975
$ws.Range("C29").Value = "'01/10/2019"
$ws.Range("D29").Value = 29.3125853853587
$ws.Range("C30").Value = "'01/10/2021"
$ws.Range("D30").Value = 30.26671294819596
$ws.Range("C31").Value = "'01/10/2023"
$ws.Range("D31").Value = 30.67049781837195
$ws.Range("C32").Value = "'01/10/2013"
$ws.Range("D32").Value = 76.49234115834315
$ws.Range("C33").Value = "'01/10/2015"
$ws.Range("D33").Value = 77.44565704352813
$ws.Range("C34").Value = "'01/10/2017"
$ws.Range("D34").Value = 78.28675727880065
$ws.Range("C35").Value = "'01/10/2019"
$ws.Range("D35").Value = 78.96286946460502
$ws.Range("C36").Value = "'01/10/2021"
$ws.Range("C37").Value = "'01/10/2023"
$ws.Range("D37").Value = 79.95564800330061
$ws.Range("C38").Value = "'01/10/2013"
$ws.Range("D38").Value = 44.18381219976435
$ws.Range("C39").Value = "'01/10/2015"
$ws.Range("D39").Value = 44.71981793093561
$ws.Range("C40").Value = "'01/10/2017"
$ws.Range("D40").Value = 43.98291719090571
$ws.Range("C41").Value = "'01/10/2019"
$ws.Range("D41").Value = 44.20076434907612
$ws.Range("C42").Value = "'01/10/2021"
$ws.Range("C43").Value = "'01/10/2023"
$ws.Range("D43").Value = 43.34118375135377
$ws.Range("C44").Value = "'01/10/2013"
$ws.Range("D44").Value = 40.66527689658298
$ws.Range("C45").Value = "'01/10/2015"
$ws.Range("D45").Value = 39.98351343117753
$ws.Range("C46").Value = "'01/10/2017"
$ws.Range("D46").Value = 37.83735889847779
$ws.Range("C47").Value = "'01/10/2019"
$ws.Range("D47").Value = 38.11928053013569
$ws.Range("C48").Value = "'01/10/2021"
$ws.Range("C49").Value = "'01/10/2023"
$ws.Range("D49").Value = 38.81831152980007
$ws.Range("C50").Value = "'01/10/2013"
$ws.Range("D50").Value = 3.518535303181365
$ws.Range("C51").Value = "'01/10/2015"
$ws.Range("D51").Value = 4.736304499758077
$ws.Range("C52").Value = "'01/10/2017"
$ws.Range("D52").Value = 6.145558292427921
$ws.Range("C53").Value = "'01/10/2019"
$ws.Range("D53").Value = 6.079730724729147
$ws.Range("C54").Value = "'01/10/2021"
$ws.Range("C55").Value = "'01/10/2023"
$ws.Range("D55").Value = 4.522872221553695
$ws.Range("C56").Value = "'01/10/2013"
$ws.Range("D56").Value = 32.30671621499139
$ws.Range("C57").Value = "'01/10/2015"
$ws.Range("D57").Value = 32.72583911259252
$ws.Range("C58").Value = "'01/10/2017"
$ws.Range("D58").Value = 34.30384008789495
$ws.Range("C59").Value = "'01/10/2019"
$ws.Range("D59").Value = 34.76210511552891
$ws.Range("C60").Value = "'01/10/2021"
$ws.Range("C61").Value = "'01/10/2023"
$ws.Range("D61").Value = 36.61446425194685
$ws.Range("C62").Value = "'01/10/2013"
$ws.Range("D62").Value = 76.29527739569005
$ws.Range("C63").Value = "'01/10/2015"
$ws.Range("D63").Value = 78.50719424460432
$ws.Range("C64").Value = "'01/10/2017"
$ws.Range("D64").Value = 78.02294792586054
$ws.Range("C65").Value = "'01/10/2019"
$ws.Range("D65").Value = 78.75108412836079
$ws.Range("C66").Value = "'01/10/2021"
$ws.Range("C67").Value = "'01/10/2023"
$ws.Range("D67").Value = 78.59840537138061
$ws.Range("C68").Value = "'01/10/2013"
$ws.Range("D68").Value = 45.89637780834479
$ws.Range("C69").Value = "'01/10/2015"
$ws.Range("D69").Value = 45.45863309352518
$ws.Range("C70").Value = "'01/10/2017"
$ws.Range("D70").Value = 46.07237422771404
$ws.Range("C71").Value = "'01/10/2019"
$ws.Range("D71").Value = 48.78577623590633
$ws.Range("C72").Value = "'01/10/2021"
$ws.Range("C73").Value = "'01/10/2023"
$ws.Range("D73").Value = 44.94334872010072
$ws.Range("C74").Value = "'01/10/2013"
$ws.Range("D74").Value = 41.86153140761119
$ws.Range("C75").Value = "'01/10/2015"
$ws.Range("D75").Value = 40.87230215827338
$ws.Range("C76").Value = "'01/10/2017"
$ws.Range("D76").Value = 39.84995586937335
$ws.Range("C77").Value = "'01/10/2019"
$ws.Range("D77").Value = 41.50043365134432
$ws.Range("C78").Value = "'01/10/2021"
$ws.Range("C79").Value = "'01/10/2023"
$ws.Range("D79").Value = 39.86571548468317
$ws.Range("C80").Value = "'01/10/2013"
$ws.Range("D80").Value = 4.034846400733608
$ws.Range("C81").Value = "'01/10/2015"
$ws.Range("D81").Value = 4.631294964028776
$ws.Range("C82").Value = "'01/10/2017"
$ws.Range("D82").Value = 6.266548984995588
$ws.Range("C83").Value = "'01/10/2019"
$ws.Range("D83").Value = 7.285342584562011
$ws.Range("C84").Value = "'01/10/2021"
$ws.Range("C85").Value = "'01/10/2023"
$ws.Range("D85").Value = 5.035669324381032
$ws.Range("C86").Value = "'01/10/2013"
$ws.Range("D86").Value = 30.39889958734526
$ws.Range("C87").Value = "'01/10/2015"
$ws.Range("D87").Value = 33.04856115107913
$ws.Range("C88").Value = "'01/10/2017"
$ws.Range("D88").Value = 31.95057369814651
$ws.Range("C89").Value = "'01/10/2019"
$ws.Range("D89").Value = 29.96530789245447
$ws.Range("C90").Value = "'01/10/2021"
$ws.Range("C91").Value = "'01/10/2023"
$ws.Range("D91").Value = 33.6550566512799
